$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ROWS values for the two remaining data rows
$ws.Range("C2").Value = 20000
$ws.Range("C3").Value = 20000

# Clear the now-removed table rows (4-9): values & contents, styles stay as-is
$ws.Range("A4:E9").ClearContents()
